# Update the "Corr/total marks" figures on the concise marksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row - right-answer mark value
$ws.Range("B11").Value = 5

# "Total" row - corrected total score and the "correct/total" summary text
$ws.Range("B12").Value = 100
$ws.Range("E12").Value = "100/140"
